$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.016.15'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.48%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.641.19'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.24%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.32%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.03'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.34%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5137'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.23%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.29%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2591'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.63%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06379'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.51%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.83'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.38%  '

# Row 11
$ws.Range('E11').Value = '  +0.19%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.665.29'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.26%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.297'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.24%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5491'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.50%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '64.69'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.71%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₅7768'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.59%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.049.88'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.31%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.35%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '199.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.86%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.471'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.36%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.995'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.14%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.121'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.02%  '

# Row 23
$ws.Range('E23').Value = '  -0.29%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.901'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.20%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '142.09'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.72%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1221'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.37%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.891'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.05%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.70'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.23%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.244'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.22%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.04874'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.11%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.303'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.01%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.245'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.21%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.542'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.06%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.384'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.78%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9175'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.58%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.595'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.04%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.5598'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.86%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.114.33'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.80%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01573'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.46%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.002'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.38%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.543'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.63%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.567'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.98%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8116'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.45%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.67'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.05%  '

# Row 45
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.782.47'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.01%  '

# Row 46
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0₈118'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.55%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4540'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.18%  '

# Row 48
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.008'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.11%  '

# Row 49
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '55.27'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.00%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05233'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.93%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.09575'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.02%  '
